# Registrar la nueva hoja "asesorias" en el libro de registros.
$wb = $excel.ActiveWorkbook

$wsDocentes = $wb.Worksheets.Item("docentes")
$wsNotif = $wb.Worksheets.Item("notificaciones")

# Copiar el formato de encabezado (fuente en negrita + borde) desde "docentes"
# y pegarlo en la nueva hoja, insertada justo antes de "notificaciones".
$wsDocentes.Range("A1:K1").Copy() | Out-Null

$wsAsesorias = $wb.Worksheets.Add($wsNotif)
$wsAsesorias.Name = "asesorias"

$wsAsesorias.Range("A1:K1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Encabezados de la nueva hoja de asesorías
$wsAsesorias.Range("A1").Value = "Estudiante"
$wsAsesorias.Range("B1").Value = "Asesor"
$wsAsesorias.Range("C1").Value = "Motivo asesoría"
$wsAsesorias.Range("D1").Value = "Fecha"
$wsAsesorias.Range("E1").Value = "Hora"

# Ajustar el ancho de las columnas principales al contenido
$wsAsesorias.Columns("A").ColumnWidth = 17.666666666666668
$wsAsesorias.Columns("B").ColumnWidth = 13.666666666666666
$wsAsesorias.Columns("C").ColumnWidth = 14.333333333333334

# Selección activa en la nueva hoja
$wsAsesorias.Range("F4").Select() | Out-Null
